$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new order-line rows (4-11). Values that look numeric (SKU codes,
# quantities, dollar amounts) get a leading apostrophe so Excel stores them
# as plain text, matching the existing rows on the sheet.

# Row 4: Capers
$ws.Range("A4").Value = "'10198"
$ws.Range("B4").Value = 'Capers'
$ws.Range("C4").Value = "'1"
$ws.Range("D4").Value = "'$62.18"
$ws.Range("E4").Value = "'$62.18"

# Row 5: PKT Hot Sauce - Cholula
$ws.Range("A5").Value = 'A0432'
$ws.Range("B5").Value = 'PKT Hot Sauce - Cholula'
$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "'$16.59"
$ws.Range("E5").Value = "'$16.59"

# Row 6: PKT Sugar
$ws.Range("A6").Value = 'P0998'
$ws.Range("B6").Value = 'PKT Sugar'
$ws.Range("C6").Value = "'1"
$ws.Range("D6").Value = "'$18.13"
$ws.Range("E6").Value = "'$18.13"

# Row 7: Poland Spring - Sport Top
$ws.Range("A7").Value = 'FV234'
$ws.Range("B7").Value = 'Poland Spring - Sport Top'
$ws.Range("C7").Value = "'5"
$ws.Range("D7").Value = "'$11.09"
$ws.Range("E7").Value = "'$55.45"

# Row 8: Simply - Orange
$ws.Range("A8").Value = 'AT104'
$ws.Range("B8").Value = 'Simply - Orange'
$ws.Range("C8").Value = "'1"
$ws.Range("D8").Value = "'$44.33"
$ws.Range("E8").Value = "'$44.33"

# Row 9: Sparkling Ice - Watermelon Strawberry
$ws.Range("A9").Value = 'B2402'
$ws.Range("B9").Value = 'Sparkling Ice - Watermelon Strawberry'
$ws.Range("C9").Value = "'1"
$ws.Range("D9").Value = "'$12.25"
$ws.Range("E9").Value = "'$12.25"

# Row 10: Sparkling Ice - Kiwi Strawberry
$ws.Range("A10").Value = "'99458"
$ws.Range("B10").Value = 'Sparkling Ice - Kiwi Strawberry'
$ws.Range("C10").Value = "'1"
$ws.Range("D10").Value = "'$11.15"
$ws.Range("E10").Value = "'$11.15"

# Row 11: Sparkling Ice - Black Raspberry
$ws.Range("A11").Value = "'99402"
$ws.Range("B11").Value = 'Sparkling Ice - Black Raspberry'
$ws.Range("C11").Value = "'1"
$ws.Range("D11").Value = "'$12.25"
$ws.Range("E11").Value = "'$12.25"
